$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRMAccuracyData")

# dmb ran one more test sample for run 0506 - append it as a new row
$ws.Range("A34").Value = 20210506
$ws.Range("B34").Value = 2225.0940000000001
$ws.Range("C34").Value = 2224.4699999999998
$ws.Range("D34").Formula = "=100*(B34-C34)/C34"
$ws.Range("E34").Value = 180
$ws.Range("F34").Value = "CRM opened 20210418"

# Scroll/select to mirror where the user ended up after entering the new row
$excel.ActiveWindow.ScrollRow = 22
$ws.Range("A35").Select()
